$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (shifts existing B:Y to C:Z)
$ws.Columns("B:B").Insert()

# New header for the inserted date column
$ws.Range("B1").Value = "10_11_2020"

# New counts for the inserted column, rows 2-11 (age groups 0-9 .. 90+)
$values = @(39, 59, 172, 270, 427, 704, 704, 990, 746, 196)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(2 + $i, 2).Value = $values[$i]
}

# Match formatting of column A (quote-prefixed text style) for the new numeric column
$ws.Range("A2:A11").Copy()
$ws.Range("B2:B11").PasteSpecial(-4122)  # xlPasteFormats

# Sum row for the new column
$ws.Range("B12").Formula = "=SUM(B2:B11)"

# Reset view: scroll back to column A, and move selection
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("D20").Select()
